# HeatData.xlsx edit:
#  - Duplicate Sheet1 (with its header row + 4 data rows) into a new
#    "Sheet2" placed after Sheet1.
#  - On Sheet1, remove the header row so the sheet holds just the four
#    numeric data rows (shifted up to rows 1-4), and shrink the columns
#    to a narrow, uniform width now that there is no header text driving
#    the "best fit" size.
#  - Leave the selection on Sheet2 covering the whole used range, and
#    move Sheet1's selection/active cell to G11, restoring Sheet1 as the
#    active sheet when done.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Build Sheet2 as a copy of the original Sheet1 (header + data) ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Sheet2 keeps the original layout/selection: whole data range selected,
# no single active cell recorded, and it is not the active tab.
$ws2.Range("A1:F5").Select()

# --- Strip the header row from Sheet1, shifting data up one row ---
$ws1.Rows.Item(1).Delete()

# Narrow, uniform column widths now that headers are gone (still marked
# as "best fit" in the saved file).
$ws1.Columns.Item(1).ColumnWidth = 3.1666666666666665
$ws1.Columns.Item(2).ColumnWidth = 4.166666666666667
$ws1.Columns.Item(3).ColumnWidth = 4.166666666666667
$ws1.Columns.Item(4).ColumnWidth = 3.1666666666666665
$ws1.Columns.Item(5).ColumnWidth = 3.1666666666666665
$ws1.Columns.Item(6).ColumnWidth = 4.166666666666667

# Sheet1 is active again, with the cursor parked at G11.
$ws1.Activate()
$ws1.Range("G11").Select()
